$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 6, shifting existing rows 6-60 down to 7-61.
$ws.Rows.Item(6).Insert()

# Populate the newly inserted row 6 with the new record.
$newRow = $ws.Rows.Item(6)
$newRow.Cells.Item(1, 1).Value = 9
$newRow.Cells.Item(1, 2).Value = "Vega Central Mapocho de Santiago"
$newRow.Cells.Item(1, 3).Value = "Metropolitana"
$newRow.Cells.Item(1, 4).Value = 44552
$newRow.Cells.Item(1, 5).Value = 13
$newRow.Cells.Item(1, 6).Value = "Fruta"
$newRow.Cells.Item(1, 7).Value = 100101
$newRow.Cells.Item(1, 8).Value = "Berries"
$newRow.Cells.Item(1, 9).Value = 100101004
$newRow.Cells.Item(1, 10).Value = "Frambuesa"
$newRow.Cells.Item(1, 11).Value = "Sin especificar"
$newRow.Cells.Item(1, 12).Value = "Primera"
$newRow.Cells.Item(1, 13).Value = 610
$newRow.Cells.Item(1, 14).Value = 7500
$newRow.Cells.Item(1, 15).Value = 8000
$newRow.Cells.Item(1, 16).Value = 7730
$newRow.Cells.Item(1, 17).Value = "`$/bandeja 2 kilos"
$newRow.Cells.Item(1, 18).Value = "Región del Maule"
$newRow.Cells.Item(1, 19).Value = 3865
$newRow.Cells.Item(1, 20).Value = 2

# Match the date-formatted style used by the other "Fecha" cells in column D.
$newRow.Cells.Item(1, 4).NumberFormat = $ws.Cells.Item(7, 4).NumberFormat
